$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "TESTADO E A FUNCIONAR"
$ws.Range("D8").Value = "TESTADO E A FUNCIONAR"
$ws.Range("D9").Value = "TESTADO E A FUNCIONAR"
$ws.Range("D10").Value = "x"
$ws.Range("D11").Value = "x"
$ws.Range("D12").Value = "x"
$ws.Range("D14").Value = "x"
$ws.Range("D17").Value = "x"
$ws.Range("D18").Value = "x"

$ws.Range("D9").Select()
